$wb = $excel.ActiveWorkbook

# Rename the "Include" sheet
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# Metadata sheet updates
$ws1 = $wb.Worksheets.Item(1)

# Version 3.8.0 -> 3.9.0
$ws1.Range("B3").Value = "3.9.0"

# Experimental value was blank -> false (force text, not an Excel boolean,
# by writing it as a formula result then converting it to a plain value)
$ws1.Range("B7").Formula = "=""false"""
$ws1.Range("B7").Copy()
$ws1.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Date updated
$ws1.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Contact rows (10-12), each carrying a different display string
$ws1.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws1.Range("B11").Value = "null (iti@ihe.net)"
$ws1.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction World -> Global (Whole world)
$ws1.Range("B13").Value = "Global (Whole world)"

Write-Output "done"
